$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D2 already holds 20.0; just make sure the value/style (s="0") is applied.
$ws.Range("D2").Value = 20.0

# Row 3: add Time (C3) and Material (D3)
$ws.Range("C3").Value = 1.0
$ws.Range("D3").Value = 20.0

# Row 4: add Time (C4) and Material (D4)
$ws.Range("C4").Value = 3.546999931335449
$ws.Range("D4").Value = 16.0

# Row 5: add Time (C5) and Material (D5)
$ws.Range("C5").Value = 3.421999931335449
$ws.Range("D5").Value = 51.0
